# Insert a new weekly record row at row 194 (Hortaliza, Macroferia Regional de
# Talca - Zanahoria), pushing the existing rows 194:239 down to 195:240.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(194).Insert()

$ws.Range("A194").Value = 5
$ws.Range("B194").Value = "Macroferia Regional de Talca"
$ws.Range("C194").Value = "Maule"
$ws.Range("D194").Value = 44543
$ws.Range("E194").Value = 7
$ws.Range("F194").Value = 100114013
$ws.Range("G194").Value = "Zanahoria"
$ws.Range("H194").Value = "Sin especificar"
$ws.Range("I194").Value = "Primera"
$ws.Range("J194").Value = 400
$ws.Range("K194").Value = 7000
$ws.Range("L194").Value = 7000
$ws.Range("M194").Value = 7000
$ws.Range("N194").Value = "`$/saco 20 kilos"
$ws.Range("O194").Value = "Provincia del Elquí"
$ws.Range("P194").Value = 350
$ws.Range("Q194").Value = 20
$ws.Range("R194").Value = "Hortaliza"
